$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 0.1645843333333333
$ws.Range("H2").Value = 0.493753
$ws.Range("I2").Value = 0.03485847193389392
$ws.Range("J2").Value = 0.03485847193389392
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 103.4275383333333
$ws.Range("N2").Value = 310.282615
$ws.Range("O2").Value = 0.2485530285127421
$ws.Range("P2").Value = 0.2485530285127421
$ws.Range("Q2").Value = 17.02255244489945
$ws.Range("R2").Value = 153.202972004095
$ws.Range("S2").Value = 0.008664178768495755
$ws.Range("T2").Value = 0.008664178768495755

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 0.1645843333333333
$ws.Range("H3").Value = 0.493753
$ws.Range("I3").Value = 0.03485847193389392
$ws.Range("J3").Value = 0.03485847193389392
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 216.130539
$ws.Range("N3").Value = 648.391617
$ws.Range("O3").Value = 0.5193964865470273
$ws.Range("P3").Value = 0.5193964865470272
$ws.Range("Q3").Value = 35.571700674289
$ws.Range("R3").Value = 320.145306068601
$ws.Range("S3").Value = 0.01810536784886266
$ws.Range("T3").Value = 0.01810536784886265

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 0.1645843333333333
$ws.Range("H4").Value = 0.493753
$ws.Range("I4").Value = 0.03485847193389392
$ws.Range("J4").Value = 0.03485847193389392
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 71.607325
$ws.Range("N4").Value = 214.821975
$ws.Range("O4").Value = 0.1720839321833696
$ws.Range("P4").Value = 0.1720839321833696
$ws.Range("Q4").Value = 11.78544384690833
$ws.Range("R4").Value = 106.068994622175
$ws.Range("S4").Value = 0.005998582920288092
$ws.Range("T4").Value = 0.005998582920288092

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 0.1645843333333333
$ws.Range("H5").Value = 0.493753
$ws.Range("I5").Value = 0.03485847193389392
$ws.Range("J5").Value = 0.03485847193389392
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 24.953198
$ws.Range("N5").Value = 74.859594
$ws.Range("O5").Value = 0.05996655275686102
$ws.Range("P5").Value = 0.05996655275686102
$ws.Range("Q5").Value = 4.106905457364666
$ws.Range("R5").Value = 36.962149116282
$ws.Range("S5").Value = 0.002090342396247409
$ws.Range("T5").Value = 0.002090342396247408

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 3.368329
$ws.Range("H6").Value = 10.104987
$ws.Range("I6").Value = 0.7134020567608963
$ws.Range("J6").Value = 0.7134020567608964
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 103.4275383333333
$ws.Range("N6").Value = 310.282615
$ws.Range("O6").Value = 0.2485530285127421
$ws.Range("P6").Value = 0.2485530285127421
$ws.Range("Q6").Value = 348.3779767667783
$ws.Range("R6").Value = 3135.401790901005
$ws.Range("S6").Value = 0.1773182417551399
$ws.Range("T6").Value = 0.1773182417551399

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 3.368329
$ws.Range("H7").Value = 10.104987
$ws.Range("I7").Value = 0.7134020567608963
$ws.Range("J7").Value = 0.7134020567608964
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 216.130539
$ws.Range("N7").Value = 648.391617
$ws.Range("O7").Value = 0.5193964865470273
$ws.Range("P7").Value = 0.5193964865470272
$ws.Range("Q7").Value = 727.9987622993309
$ws.Range("R7").Value = 6551.988860693978
$ws.Range("S7").Value = 0.3705385217770325
$ws.Range("T7").Value = 0.3705385217770324

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 3.368329
$ws.Range("H8").Value = 10.104987
$ws.Range("I8").Value = 0.7134020567608963
$ws.Range("J8").Value = 0.7134020567608964
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 71.607325
$ws.Range("N8").Value = 214.821975
$ws.Range("O8").Value = 0.1720839321833696
$ws.Range("P8").Value = 0.1720839321833696
$ws.Range("Q8").Value = 241.197029409925
$ws.Range("R8").Value = 2170.773264689325
$ws.Range("S8").Value = 0.1227650311551184
$ws.Range("T8").Value = 0.1227650311551185

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 3.368329
$ws.Range("H9").Value = 10.104987
$ws.Range("I9").Value = 0.7134020567608963
$ws.Range("J9").Value = 0.7134020567608964
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 24.953198
$ws.Range("N9").Value = 74.859594
$ws.Range("O9").Value = 0.05996655275686102
$ws.Range("P9").Value = 0.05996655275686102
$ws.Range("Q9").Value = 84.05058046614199
$ws.Range("R9").Value = 756.455224195278
$ws.Range("S9").Value = 0.04278026207360545
$ws.Range("T9").Value = 0.04278026207360545

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 1.188588333333333
$ws.Range("H10").Value = 3.565765
$ws.Range("I10").Value = 0.2517394713052097
$ws.Range("J10").Value = 0.2517394713052098
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 103.4275383333333
$ws.Range("N10").Value = 310.282615
$ws.Range("O10").Value = 0.2485530285127421
$ws.Range("P10").Value = 0.2485530285127421
$ws.Range("Q10").Value = 122.9327654083861
$ws.Range("R10").Value = 1106.394888675475
$ws.Range("S10").Value = 0.06257060798910642
$ws.Range("T10").Value = 0.06257060798910642

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 1.188588333333333
$ws.Range("H11").Value = 3.565765
$ws.Range("I11").Value = 0.2517394713052097
$ws.Range("J11").Value = 0.2517394713052098
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 216.130539
$ws.Range("N11").Value = 648.391617
$ws.Range("O11").Value = 0.5193964865470273
$ws.Range("P11").Value = 0.5193964865470272
$ws.Range("Q11").Value = 256.890237132445
$ws.Range("R11").Value = 2312.012134192005
$ws.Range("S11").Value = 0.1307525969211321
$ws.Range("T11").Value = 0.1307525969211321

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 1.188588333333333
$ws.Range("H12").Value = 3.565765
$ws.Range("I12").Value = 0.2517394713052097
$ws.Range("J12").Value = 0.2517394713052098
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 71.607325
$ws.Range("N12").Value = 214.821975
$ws.Range("O12").Value = 0.1720839321833696
$ws.Range("P12").Value = 0.1720839321833696
$ws.Range("Q12").Value = 85.11163107620834
$ws.Range("R12").Value = 766.004679685875
$ws.Range("S12").Value = 0.04332031810796302
$ws.Range("T12").Value = 0.04332031810796303

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 1.188588333333333
$ws.Range("H13").Value = 3.565765
$ws.Range("I13").Value = 0.2517394713052097
$ws.Range("J13").Value = 0.2517394713052098
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 24.953198
$ws.Range("N13").Value = 74.859594
$ws.Range("O13").Value = 0.05996655275686102
$ws.Range("P13").Value = 0.05996655275686102
$ws.Range("Q13").Value = 29.65908002215667
$ws.Range("R13").Value = 266.93172019941
$ws.Range("S13").Value = 0.01509594828700816
$ws.Range("T13").Value = 0.01509594828700816
